$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Trey Murphy III", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
